$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> D (price) and E (volume) new text values.
# $null means "no change for that column in this row".
$updates = @(
    @{ Row = 2;  D = "66.850.80"; E = "  -2.14%  " },
    @{ Row = 3;  D = "2.621.02";  E = "  -3.48%  " },
    @{ Row = 4;  D = $null;       E = "  -0.05%  " },
    @{ Row = 5;  D = "589.04";    E = "  -2.91%  " },
    @{ Row = 6;  D = "164.67";    E = "  -1.24%  " },
    @{ Row = 7;  D = $null;       E = "  -0.01%  " },
    @{ Row = 8;  D = "0.533";     E = "  -3.79%  " },
    @{ Row = 9;  D = "2.621.16";  E = "  -3.48%  " },
    @{ Row = 10; D = $null;       E = "  -1.73%  " },
    @{ Row = 11; D = "0.160";     E = "  +1.26%  " },
    @{ Row = 12; D = "0.360";     E = "  -1.14%  " },
    @{ Row = 13; D = $null;       E = "  -1.71%  " },
    @{ Row = 14; D = "27.45";     E = "  -3.60%  " },
    @{ Row = 15; D = "3.117.25";  E = "  -3.19%  " },
    @{ Row = 16; D = $null;       E = "  -3.63%  " },
    @{ Row = 17; D = "66.738.51"; E = "  -2.28%  " },
    @{ Row = 18; D = "2.619.85";  E = "  -3.74%  " },
    @{ Row = 19; D = "11.91";     E = "  +0.42%  " },
    @{ Row = 20; D = "8.03";      E = "  +5.32%  " },
    @{ Row = 21; D = "357.18";    E = "  -3.86%  " },
    @{ Row = 22; D = $null;       E = "  -4.05%  " },
    @{ Row = 23; D = "4.69";      E = "  -5.55%  " },
    @{ Row = 24; D = "10.95";     E = "  +7.70%  " },
    @{ Row = 25; D = "1.95";      E = "  -6.07%  " },
    @{ Row = 26; D = "1.00";      E = "  -0.02%  " },
    @{ Row = 27; D = "70.63";     E = "  -3.31%  " },
    @{ Row = 28; D = $null;       E = "  -3.48%  " },
    @{ Row = 29; D = "1.00";      E = "  +0.08%  " },
    @{ Row = 30; D = $null;       E = "  -3.21%  " },
    @{ Row = 31; D = "549.61";    E = "  -4.85%  " },
    @{ Row = 32; D = $null;       E = "  -3.41%  " },
    @{ Row = 33; D = $null;       E = "  -4.49%  " },
    @{ Row = 34; D = $null;       E = "  -5.14%  " },
    @{ Row = 35; D = $null;       E = "  +1.29%  " },
    @{ Row = 36; D = $null;       E = "  -0.07%  " },
    @{ Row = 37; D = "1.51";      E = "  -4.93%  " },
    @{ Row = 38; D = "157.18";    E = "  -3.04%  " },
    @{ Row = 39; D = $null;       E = "  -3.80%  " },
    @{ Row = 40; D = $null;       E = "  -3.18%  " },
    @{ Row = 41; D = $null;       E = "  -3.85%  " },
    @{ Row = 42; D = $null;       E = "  -4.50%  " },
    @{ Row = 43; D = $null;       E = "  -0.59%  " },
    @{ Row = 44; D = $null;       E = "  -0.02%  " },
    @{ Row = 45; D = $null;       E = "  -5.39%  " },
    @{ Row = 46; D = "40.17";     E = "  -1.72%  " },
    @{ Row = 47; D = $null;       E = "  -4.99%  " },
    @{ Row = 48; D = "0.585";     E = "  -1.74%  " },
    @{ Row = 49; D = "151.81";    E = "  -1.89%  " },
    @{ Row = 50; D = "3.80";      E = $null },
    @{ Row = 51; D = $null;       E = "  -3.19%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $cell = $ws.Cells.Item($u.Row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.Style = "Normal"
    }
}
